$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = "64d97a2435ca528474a9ee1b62a5969d"
$ws.Range("B74").Value = "7ab7fef2fd4db72bbdb0720aafcbc718"
$ws.Range("B89").Value = "677808ed7f974be62cdfb69b4daed013"
$ws.Range("B99").Value = "7295799e6758bfbfe9f01c6adf1aca08"
$ws.Range("B110").Value = "1bd5e3b761a52acf1e20b0c69324b0d1"
$ws.Range("B136").Value = "5e3fe43d9be5b777179b6c69eea2d63f"
$ws.Range("B159").Value = "17e6f09fd8ea8a8972bc475df817080f"
$ws.Range("B161").Value = "10f1715cd7ab53d5a3f38c26ac1e512f"
$ws.Range("B168").Value = "b59d55c420b531bf2814747715b21456"
$ws.Range("B169").Value = "6afcb86346c0f16cac73003425cae14d"
$ws.Range("B227").Value = "366679d9cd102f7c634ebffd2d642faa"
$ws.Range("B232").Value = "2ad3ae0d1889ca9238638c3c5377ba7a"
$ws.Range("B278").Value = "6ca2b727497da9da297e10d0e74f11fc"
$ws.Range("B281").Value = "91d6cecafdef3ad37838abc58fd1f3c8"
$ws.Range("B302").Value = "d263c9cd625e0cc36308d3fec4350e23"
$ws.Range("B339").Value = "1e506b1f2a033ed20095cbdd53afc20a"
$ws.Range("B345").Value = "1d0565d3900a06151050ed3f0730ef7c"
$ws.Range("B419").Value = "2ee5add6736bc97726d8045230c25adb"
$ws.Range("B460").Value = "0cd8625297c32aba25b0f61545f1b53e"
$ws.Range("B478").Value = "0e421a028fe726870a018a31b7132a98"
$ws.Range("B500").Value = "59328d6fbee2ac587678815c09af1874"
$ws.Range("B501").Value = "2f3dfc70d7f041da9765e62f76ca913a"
$ws.Range("B517").Value = "4411e56c2ff7e6ec8787d8f6be166e8b"
$ws.Range("B550").Value = "345984d1f1a72d556b2fb2538b0e94aa"
$ws.Range("B616").Value = "cf51451dd6f5b3073cd680b0a9c8f098"
$ws.Range("B627").Value = "cd0f810a0814b71df06adc86d49f9165"
$ws.Range("B680").Value = "dfc9b3ba408aa959d34138ce25d08e59"
$ws.Range("B685").Value = "5bc27490b7119c501eed5d24ed5b0700"
$ws.Range("B700").Value = "c1be0d083ce0ad19eb1f14e63dd5771f"
$ws.Range("B703").Value = "19cfb9675ed60fea946f78fdbeb00be0"
$ws.Range("B768").Value = "fa3438559eb36bcd278952eeb9ffd616"
$ws.Range("B816").Value = "e156ff61a68c1b859d559b0ba2bd01c0"
$ws.Range("B819").Value = "f918429f8f38492013789bfd11f54108"
$ws.Range("B825").Value = "74f20965bca711405d4b5008fd88b85c"
$ws.Range("B827").Value = "7c0d8b2c888ea89da57dac14fe891e28"
$ws.Range("B830").Value = "39131b3cfdad3487567b097fc174ea20"
$ws.Range("B835").Value = "6c0c01f6b02ef111a430a37b418b5556"
